$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PromptResults")

$ws.Range("C2").Value = 18.47
$ws.Range("C3").Value = 16.26
$ws.Range("C4").Value = 44.89
$ws.Range("C5").Value = 17.55
$ws.Range("C6").Value = 16.22
$ws.Range("C7").Value = 22.33
$ws.Range("C8").Value = 39.14
$ws.Range("C9").Value = 38.13
